$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 19
$ws.Range("I2").Value = 500
$ws.Range("J2").Value = 50000

# Row 3
$ws.Range("C3").Value = 19
$ws.Range("I3").Value = 200
$ws.Range("J3").Value = 20000

# Row 4
$ws.Range("C4").Value = 19
$ws.Range("I4").Value = 200
$ws.Range("J4").Value = 20000

# Update the active selection/cursor position to J12
$ws.Range("J12").Select()
